# Weekly update: shift existing weekly price rows down by one week (one
# Primera/Segunda pair = 2 rows), insert a brand-new week at the top of the
# data block (rows 8-9, date 44812) and push the oldest week that falls off
# the bottom of the original range into two brand-new rows (166-167).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Snapshot the current data block (rows 8-165, all columns A:R) BEFORE
#    anything is overwritten.
$snapshot = $ws.Range("A8:R165").Value2

# 2. Write the snapshot back starting two rows lower (rows 10-167). Because
#    we operate on an in-memory snapshot rather than re-reading cells as we
#    go, this performs a clean "shift down by one pair" for the whole block
#    in a single assignment, including the two brand-new trailing rows
#    (166-167) that now hold what used to be the oldest week (rows 164-165).
$ws.Range("A10:R167").Value2 = $snapshot

# 3. The two brand-new rows (166-167) fall outside the sheet's original
#    dimensions, so they don't automatically inherit the date-column number
#    format used throughout column D. Copy it explicitly.
$ws.Range("D166").NumberFormat = $ws.Range("D164").NumberFormat
$ws.Range("D167").NumberFormat = $ws.Range("D165").NumberFormat

# 4. Populate the new week inserted at rows 8-9 (date 44812). Only the
#    columns that actually change for a new week (Fecha + the four price
#    columns) need to be touched; everything else keeps the values that were
#    already sitting in rows 8-9 before the shift (Mercado, Categoria,
#    Calidad, Volumen, Unidad, Origen, etc. all repeat unchanged week over
#    week for this product/market combination).

# Row 8 = "Primera" quality week of 2022-09-08 (serial date 44812)
$ws.Range("D8").Value2 = 44812
$ws.Range("K8").Value2 = 700
$ws.Range("L8").Value2 = 800
$ws.Range("M8").Value2 = 750
$ws.Range("P8").Value2 = 750

# Row 9 = "Segunda" quality week of 2022-09-08 (serial date 44812)
$ws.Range("D9").Value2 = 44812
$ws.Range("K9").Value2 = 600
$ws.Range("L9").Value2 = 600
$ws.Range("M9").Value2 = 600
$ws.Range("P9").Value2 = 600
